# Apply the changes described by the diff to sets.xlsx

$wb = $excel.ActiveWorkbook

# --- Workbook-level properties ---

# Rename sheet "_set_a_double" (sheetId 6) to "_set_a2"
$wsADouble = $wb.Worksheets.Item("_set_a_double")
$wsADouble.Name = "_set_a2"

# --- Sheet "_set_a" (sheetId 5): no content change, only loses tabSelected
# below (handled automatically once _set_a2 is Activate()'d) ---

# --- Sheet "_set_a2" (formerly "_set_a_double") ---
$wsA2 = $wb.Worksheets.Item("_set_a2")
$wsA2.Range("A1").Value = "a2_names"

# Make "_set_a2" the active / tab-selected sheet and set its selection
$wsA2.Activate()
$wsA2.Range("G21").Select()
